# Auto-generated edit script: updates currentAveragePrice / Leve price / profit
# columns (H-N) across multiple worksheets per the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value2 = 1429365.2
$ws.Range("I33").Value2 = 3333482
$ws.Range("K33").Value2 = 3333482
$ws.Range("M33").Value2 = -3333253
$ws.Range("H40").Value2 = 7422.875
$ws.Range("J40").Value2 = 19999
$ws.Range("L40").Value2 = 19999
$ws.Range("N40").Value2 = -20349
$ws.Range("H41").Value2 = 264.42856
$ws.Range("I41").Value2 = 338.6
$ws.Range("J41").Value2 = 79
$ws.Range("K41").Value2 = 338.6
$ws.Range("L41").Value2 = 79
$ws.Range("M41").Value2 = 101.4
$ws.Range("N41").Value2 = -959
$ws.Range("H86").Value2 = 946.75
$ws.Range("I86").Value2 = 931.6667
$ws.Range("J86").Value2 = 992
$ws.Range("K86").Value2 = 931.6667
$ws.Range("L86").Value2 = 992
$ws.Range("M86").Value2 = 191.3333
$ws.Range("N86").Value2 = -3238
$ws.Range("H89").Value2 = 946.75
$ws.Range("I89").Value2 = 931.6667
$ws.Range("J89").Value2 = 992
$ws.Range("K89").Value2 = 4658.3335
$ws.Range("L89").Value2 = 992
$ws.Range("M89").Value2 = 957.6665000000003
$ws.Range("N89").Value2 = -16192
$ws.Range("H106").Value2 = 11321.138
$ws.Range("I106").Value2 = 12360.385
$ws.Range("K106").Value2 = 12360.385
$ws.Range("M106").Value2 = -11729.385
$ws.Range("H137").Value2 = 23812438
$ws.Range("I137").Value2 = 41668904
$ws.Range("K137").Value2 = 125006712
$ws.Range("M137").Value2 = -125004162
$ws.Range("H138").Value2 = 2259.3076
$ws.Range("I138").Value2 = 1263.4
$ws.Range("J138").Value2 = 2881.75
$ws.Range("K138").Value2 = 3790.2
$ws.Range("L138").Value2 = 8645.25
$ws.Range("M138").Value2 = 1349.8
$ws.Range("N138").Value2 = -18925.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value2 = 9664
$ws.Range("I32").Value2 = 7191.5835
$ws.Range("J32").Value2 = 15597.8
$ws.Range("K32").Value2 = 7191.5835
$ws.Range("L32").Value2 = 15597.8
$ws.Range("M32").Value2 = -6904.5835
$ws.Range("N32").Value2 = -16171.8
$ws.Range("H45").Value2 = 2106
$ws.Range("I45").Value2 = 2127.2
$ws.Range("K45").Value2 = 2127.2
$ws.Range("M45").Value2 = -1750.2
$ws.Range("H74").Value2 = 1577.18
$ws.Range("I74").Value2 = 1510.65
$ws.Range("K74").Value2 = 1510.65
$ws.Range("M74").Value2 = -636.6500000000001
$ws.Range("H77").Value2 = 1577.18
$ws.Range("I77").Value2 = 1510.65
$ws.Range("K77").Value2 = 7553.25
$ws.Range("M77").Value2 = -3185.25
$ws.Range("H132").Value2 = 7598.5396
$ws.Range("I132").Value2 = 7633.017
$ws.Range("J132").Value2 = 7198.6
$ws.Range("K132").Value2 = 22899.051
$ws.Range("L132").Value2 = 21595.8
$ws.Range("M132").Value2 = -20369.051
$ws.Range("N132").Value2 = -26655.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value2 = 1328.625
$ws.Range("I5").Value2 = 1435.8
$ws.Range("K5").Value2 = 1435.8
$ws.Range("M5").Value2 = -1322.8
$ws.Range("H20").Value2 = 910.6667
$ws.Range("I20").Value2 = 954.3570999999999
$ws.Range("K20").Value2 = 954.3570999999999
$ws.Range("M20").Value2 = -707.3570999999999
$ws.Range("H81").Value2 = 41609.832
$ws.Range("J81").Value2 = 41609.832
$ws.Range("L81").Value2 = 41609.832
$ws.Range("N81").Value2 = -43731.832
$ws.Range("H84").Value2 = 41609.832
$ws.Range("J84").Value2 = 41609.832
$ws.Range("L84").Value2 = 124829.496
$ws.Range("N84").Value2 = -135437.496
$ws.Range("H99").Value2 = 58831770
$ws.Range("I99").Value2 = 58831770
$ws.Range("K99").Value2 = 58831770
$ws.Range("M99").Value2 = -58830272
$ws.Range("H134").Value2 = 11815.787
$ws.Range("I134").Value2 = 10492.463
$ws.Range("J134").Value2 = 18636
$ws.Range("K134").Value2 = 31477.389
$ws.Range("L134").Value2 = 55908
$ws.Range("M134").Value2 = -28942.389
$ws.Range("N134").Value2 = -60978

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value2 = 906.75
$ws.Range("I12").Value2 = 405
$ws.Range("J12").Value2 = 1074
$ws.Range("K12").Value2 = 405
$ws.Range("L12").Value2 = 1074
$ws.Range("M12").Value2 = -235
$ws.Range("N12").Value2 = -1414
$ws.Range("H22").Value2 = 3961.2727
$ws.Range("I22").Value2 = 3999.375
$ws.Range("K22").Value2 = 3999.375
$ws.Range("M22").Value2 = -3649.375
$ws.Range("H39").Value2 = 2849.4285
$ws.Range("I39").Value2 = 3158
$ws.Range("K39").Value2 = 3158
$ws.Range("M39").Value2 = -2767
$ws.Range("H49").Value2 = 2849.4285
$ws.Range("I49").Value2 = 3158
$ws.Range("K49").Value2 = 3158
$ws.Range("M49").Value2 = -2976

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value2 = 2344.7917
$ws.Range("J34").Value2 = 3699.2856
$ws.Range("L34").Value2 = 11097.8568
$ws.Range("N34").Value2 = -11265.8568
$ws.Range("H38").Value2 = 112.14286
$ws.Range("J38").Value2 = 79.111115
$ws.Range("L38").Value2 = 237.333345
$ws.Range("N38").Value2 = -931.333345
$ws.Range("H113").Value2 = 2255.2666
$ws.Range("I113").Value2 = 2757.6
$ws.Range("J113").Value2 = 2004.1
$ws.Range("K113").Value2 = 8272.799999999999
$ws.Range("L113").Value2 = 6012.299999999999
$ws.Range("M113").Value2 = -6102.799999999999
$ws.Range("N113").Value2 = -10352.3
$ws.Range("H114").Value2 = 17464.75
$ws.Range("I114").Value2 = 2875
$ws.Range("J114").Value2 = 24096.455
$ws.Range("K114").Value2 = 8625
$ws.Range("L114").Value2 = 72289.36500000001
$ws.Range("M114").Value2 = -5371
$ws.Range("N114").Value2 = -78797.36500000001
$ws.Range("H122").Value2 = 632.0769
$ws.Range("J122").Value2 = 572.125
$ws.Range("L122").Value2 = 5149.125
$ws.Range("N122").Value2 = -10049.125
$ws.Range("H124").Value2 = 8158.4287
$ws.Range("I124").Value2 = 1027.75
$ws.Range("J124").Value2 = 17666
$ws.Range("K124").Value2 = 3083.25
$ws.Range("L124").Value2 = 52998
$ws.Range("M124").Value2 = 1826.75
$ws.Range("N124").Value2 = -62818
$ws.Range("H132").Value2 = 3249.1428
$ws.Range("I132").Value2 = 2929.8
$ws.Range("K132").Value2 = 26368.2
$ws.Range("M132").Value2 = -23838.2
$ws.Range("H139").Value2 = 5749.6
$ws.Range("I139").Value2 = 4083.3333
$ws.Range("K139").Value2 = 12249.9999
$ws.Range("M139").Value2 = -7109.999899999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value2 = 1119.75
$ws.Range("I31").Value2 = 1119.75
$ws.Range("K31").Value2 = 1119.75
$ws.Range("M31").Value2 = -827.75
$ws.Range("H37").Value2 = 1119.75
$ws.Range("I37").Value2 = 1119.75
$ws.Range("K37").Value2 = 1119.75
$ws.Range("M37").Value2 = -842.75
$ws.Range("H43").Value2 = 53994.5
$ws.Range("I43").Value2 = 28000
$ws.Range("K43").Value2 = 28000
$ws.Range("M43").Value2 = -27849
$ws.Range("H80").Value2 = 2578.389
$ws.Range("I80").Value2 = 2644.2856
$ws.Range("J80").Value2 = 2347.75
$ws.Range("K80").Value2 = 2644.2856
$ws.Range("L80").Value2 = 2347.75
$ws.Range("M80").Value2 = -1646.2856
$ws.Range("N80").Value2 = -4343.75
$ws.Range("H83").Value2 = 2578.389
$ws.Range("I83").Value2 = 2644.2856
$ws.Range("J83").Value2 = 2347.75
$ws.Range("K83").Value2 = 13221.428
$ws.Range("L83").Value2 = 11738.75
$ws.Range("M83").Value2 = -8229.428
$ws.Range("N83").Value2 = -21722.75
$ws.Range("H126").Value2 = 9755383
$ws.Range("I126").Value2 = 3991
$ws.Range("K126").Value2 = 11973
$ws.Range("M126").Value2 = -9503

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H19").Value2 = 56000
$ws.Range("J19").Value2 = 56000
$ws.Range("L19").Value2 = 56000
$ws.Range("N19").Value2 = -56340
$ws.Range("H82").Value2 = 2217.4211
$ws.Range("I82").Value2 = 901.4167
$ws.Range("K82").Value2 = 901.4167
$ws.Range("M82").Value2 = -540.4167
$ws.Range("H85").Value2 = 2217.4211
$ws.Range("I85").Value2 = 901.4167
$ws.Range("K85").Value2 = 901.4167
$ws.Range("M85").Value2 = 346.5833

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value2 = 7938.375
$ws.Range("I17").Value2 = 8715.286
$ws.Range("K17").Value2 = 8715.286
$ws.Range("M17").Value2 = -8543.286
$ws.Range("H54").Value2 = 23019.4
$ws.Range("I54").Value2 = 0
$ws.Range("J54").Value2 = 23019.4
$ws.Range("K54").Value2 = 0
$ws.Range("L54").Value2 = 23019.4
$ws.Range("M54").ClearContents()
$ws.Range("N54").Value2 = -24059.4
$ws.Range("H81").Value2 = 2870.8572
$ws.Range("I81").Value2 = 2870.8572
$ws.Range("J81").Value2 = 0
$ws.Range("K81").Value2 = 5741.7144
$ws.Range("L81").Value2 = 0
$ws.Range("M81").Value2 = -4680.7144
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value2 = 2870.8572
$ws.Range("I84").Value2 = 2870.8572
$ws.Range("J84").Value2 = 0
$ws.Range("K84").Value2 = 28708.572
$ws.Range("L84").Value2 = 0
$ws.Range("M84").Value2 = -23404.572
$ws.Range("N84").ClearContents()
